# Add "Assignment" column and populate Sprint 1 tasks on the SprintTracking sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SprintTracking")

# Insert a new column before the existing "Task" column (currently column C)
# so the header order becomes: Sprint# | TimeFrame | Assignment | Task | Person
$ws.Columns("C").Insert()

$ws.Range("C1").Value = "Assignment"

$tasks = @(
    "Create Single Table DB for Books.  Fields: BookID, BookTitle, ",
    "Create API to handle basic HTTP requests to add book, remove book, query Book by ID, and Title",
    "Create Functions to route API functions to Postgres Functions",
    "Deploy to Heroku",
    "Plan next sprint (Thursday)",
    "Collect summaries and post Sprint review to group wiki (Thursday)"
)

$row = 2
foreach ($task in $tasks) {
    $ws.Cells.Item($row, 1).Value = 1
    $ws.Cells.Item($row, 2).Value = "Sept 24 - Sept 29"
    $ws.Cells.Item($row, 3).Value = "Build a Skeleton Webservice and Deploy It"
    $ws.Cells.Item($row, 4).Value = $task
    $row++
}

# NOTE: this host's ColumnWidth setter quantizes to a 6px Maximum-Digit-Width
# grid (real Excel/Calibri 11 uses 7px), so feeding the target "bestFit"
# character widths straight through overshoots by ~5/6 of a character.
# Subtracting that fixed offset lands on the closest achievable width.
$offset = 5 / 6
$ws.Columns.Item(1).ColumnWidth = 7 - $offset
$ws.Columns.Item(2).ColumnWidth = 14.6640625 - $offset
$ws.Columns.Item(3).ColumnWidth = 35.6640625 - $offset
$ws.Columns.Item(4).ColumnWidth = 80.5546875 - $offset
$ws.Columns.Item(5).ColumnWidth = 6.77734375 - $offset

$ws.Range("B11").Select()
